$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://down-id.img.susercontent.com/file/9fba4be140a535f2ad57c8829fe36a80"
$ws.Range("B2").Value = "Miniso Official Boneka Small Penguin Plush Toy Boneka Lucu mainan anak boneka lucu lembut boneka gemoy boneka import boneka anak Hadiah Ulang Tahun Kado anak Kado untuk cewek Hadiah untuk cowok Kado Natal"
$ws.Range("C2").Value = "Rp47.900"

$desc = "Miniso Official Product`nSmall Penguin Plush Toy`nDeskripsi:`n1. 2 Variasi : Abu-abu & Biru`n2. Bahan :`nFabric: Cover: 93% Polyester, 7% Spandex`nFilling: 100% Polyester`n3. Ukuran : 28/30/43 cm`n4. Efek dekoratif yang baik, nyaman dan indah.`n5. Nyaman Dipeluk dan di jadikan sebagai bantal"
$ws.Range("D2").Value = $desc

$ws.Range("E2").Value = "KAB. KARAWANG"
